$d = $word.ActiveDocument

# Locate the paragraph "FECHA DE REALIZACIÓN DEL PROTOCOLO: «=current_date»"
# The word "PROTOCOLO" needs to become "CONSENTIMIENTO", and (per the target
# OOXML) the run that used to hold "REALIZACIÓN DEL PROTOCOLO" ends up split
# into "REALIZACIÓN DEL " + a brand-new "CONSENTIMIENTO" run, while the
# neighboring "FECHA DE " / ":" / " " runs stay distinct runs (not merged
# into the edited run) exactly as they were before the edit.

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "FECHA DE REALIZACIÓN DEL PROTOCOLO"
$find.Forward = $true
$find.Wrap = 0
$found = $find.Execute()

$pStart = $find.Parent.Start
$pEnd   = $find.Parent.End

# Offsets (characters) within the matched range:
#   0 .. 9   -> "FECHA DE "
#   9 .. 25  -> "REALIZACIÓN DEL "
#   25 .. 34 -> "PROTOCOLO"
$protocolStart = $pStart + 25

# 1) Swap "PROTOCOLO" for "CONSENTIMIENTO".
$rProtocolo = $d.Range($protocolStart, $pEnd)
$rProtocolo.Text = "CONSENTIMIENTO"

$consentLen = 14  # len("CONSENTIMIENTO")
$rConsentimiento = $d.Range($protocolStart, $protocolStart + $consentLen)

# Re-apply bold explicitly so this newly-typed text keeps its own run
# instead of being silently re-merged with the preceding "REALIZACIÓN DEL ".
$rConsentimiento.Font.Bold = $false
$rConsentimiento.Font.Bold = $true

# 2) Keep "FECHA DE " as its own run (don't let it re-absorb "REALIZACIÓN DEL ").
$rFecha = $d.Range($pStart, $pStart + 9)
$rFecha.Font.Bold = $false
$rFecha.Font.Bold = $true

# 3) Keep ":" as its own run, separate from the following " ".
$colonStart = $protocolStart + $consentLen
$rColon = $d.Range($colonStart, $colonStart + 1)
$rColon.Font.Bold = $false
$rColon.Font.Bold = $true
